$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet has separate "First Name *" (col G) and "Last Name *" (col H)
# columns that need to become a single "Full Name *" column (first + " " +
# last), with the now-redundant Last Name column removed entirely (so every
# column after it shifts one to the left).
# ---------------------------------------------------------------------------

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
$deletedCol = 8                                                # column H

# --- Preserve the existing hyperlinks (the engine does not auto-shift their
#     cell references when a column is deleted), so they can be restored -
#     pointing at the correct, possibly shifted, cell - after the edit.
#     Every existing hyperlink here targets "mailto:" + the cell's own
#     display text.
$hlCount = $ws.Hyperlinks.Count
$hlRows = @()
$hlCols = @()
$hlTexts = @()
for ($i = 1; $i -le $hlCount; $i++) {
    $hl = $ws.Hyperlinks.Item($i)
    $hlRows += $hl.Range.Row
    $hlCols += $hl.Range.Column
    $hlTexts += $hl.Range.Value2
}
if ($hlCount -gt 0) {
    $ws.Cells.Hyperlinks.Delete()
}

# --- Header row: "First Name *" / "Last Name *" -> "Full Name *"
$ws.Cells.Item(1, 7).Value = "Full Name *"

# --- Data rows: concatenate first + last name with a space in between
for ($r = 2; $r -le $lastRow; $r++) {
    $first = $ws.Cells.Item($r, 7).Value2
    $last = $ws.Cells.Item($r, $deletedCol).Value2
    $ws.Cells.Item($r, 7).Value = ($first + " " + $last)
}

# --- Remove the Last Name column entirely; everything to the right shifts
#     left by one column.
$ws.Columns($deletedCol).Delete()

# --- Re-create the hyperlinks at their (possibly shifted) location, then
#     reapply the plain "Hyperlink" cell style so re-adding doesn't leave
#     the cell on a freshly minted duplicate style.
for ($i = 0; $i -lt $hlRows.Count; $i++) {
    $col = $hlCols[$i]
    if ($col -gt $deletedCol) {
        $col = $col - 1
    }
    $target = $ws.Cells.Item($hlRows[$i], $col)
    $ws.Hyperlinks.Add($target, ("mailto:" + $hlTexts[$i]))
    $target.Style = "Hyperlink"
}

# --- Update the active selection to match the post-edit state (column H now
#     holds what used to be column I - the Address column).
$ws.Range("H1:H1048576").Select()
